$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 447.5
$ws.Range("H41").Value = 445.63635
$ws.Range("I41").Value = 100
$ws.Range("J41").Value = 480.2
$ws.Range("K41").Value = 100
$ws.Range("L41").Value = 480.2
$ws.Range("M41").Value = 340
$ws.Range("N41").Value = -1360.2
$ws.Range("H43").Value = 11000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 11000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 11000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -11138
$ws.Range("H75").Value = 30866.666
$ws.Range("J75").Value = 30866.666
$ws.Range("L75").Value = 30866.666
$ws.Range("N75").Value = -32738.666
$ws.Range("H78").Value = 30866.666
$ws.Range("J78").Value = 30866.666
$ws.Range("L78").Value = 92599.99800000001
$ws.Range("N78").Value = -101959.998
$ws.Range("H112").Value = 507346.62
$ws.Range("I112").Value = 685
$ws.Range("J112").Value = 535494.5
$ws.Range("K112").Value = 2055
$ws.Range("L112").Value = 1606483.5
$ws.Range("M112").Value = -947
$ws.Range("N112").Value = -1608699.5
$ws.Range("H116").Value = 351890.38
$ws.Range("I116").Value = 1003201.5
$ws.Range("J116").Value = 9095.053
$ws.Range("K116").Value = 1003201.5
$ws.Range("L116").Value = 9095.053
$ws.Range("M116").Value = -999759.5
$ws.Range("N116").Value = -15979.053
$ws.Range("H118").Value = 430.125
$ws.Range("I118").Value = 353
$ws.Range("J118").Value = 970
$ws.Range("K118").Value = 1059
$ws.Range("L118").Value = 2910
$ws.Range("M118").Value = 598
$ws.Range("N118").Value = -6224
$ws.Range("H137").Value = 3799.0833
$ws.Range("I137").Value = 2958
$ws.Range("K137").Value = 8874
$ws.Range("M137").Value = -6324

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1474.3334
$ws.Range("I2").Value = 1533.625
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1533.625
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -1420.625
$ws.Range("N2").Value = -1226
$ws.Range("H80").Value = 37216.168
$ws.Range("J80").Value = 37216.168
$ws.Range("L80").Value = 37216.168
$ws.Range("N80").Value = -39212.168
$ws.Range("H83").Value = 37216.168
$ws.Range("J83").Value = 37216.168
$ws.Range("L83").Value = 111648.504
$ws.Range("N83").Value = -121632.504
$ws.Range("H88").Value = 6670956.5
$ws.Range("I88").Value = 9527338
$ws.Range("J88").Value = 6066.6665
$ws.Range("K88").Value = 9527338
$ws.Range("L88").Value = 6066.6665
$ws.Range("M88").Value = -9526932
$ws.Range("N88").Value = -6878.6665
$ws.Range("H91").Value = 6670956.5
$ws.Range("I91").Value = 9527338
$ws.Range("J91").Value = 6066.6665
$ws.Range("K91").Value = 9527338
$ws.Range("L91").Value = 6066.6665
$ws.Range("M91").Value = -9525934
$ws.Range("N91").Value = -8874.666499999999
$ws.Range("H116").Value = 1474.3334
$ws.Range("I116").Value = 1533.625
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 1533.625
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = 760.375
$ws.Range("N116").Value = -5588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1474.3334
$ws.Range("I3").Value = 1533.625
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 1533.625
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -1419.625
$ws.Range("N3").Value = -1228
$ws.Range("H86").Value = 2060.8635
$ws.Range("I86").Value = 1852.4546
$ws.Range("K86").Value = 1852.4546
$ws.Range("M86").Value = -729.4546
$ws.Range("H89").Value = 2060.8635
$ws.Range("I89").Value = 1852.4546
$ws.Range("K89").Value = 9262.273000000001
$ws.Range("M89").Value = -3646.273000000001
$ws.Range("H99").Value = 1856.0667
$ws.Range("I99").Value = 1048.8889
$ws.Range("J99").Value = 3066.8333
$ws.Range("K99").Value = 1048.8889
$ws.Range("L99").Value = 3066.8333
$ws.Range("M99").Value = 449.1111000000001
$ws.Range("N99").Value = -6062.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3469.4
$ws.Range("I31").Value = 1185.5454
$ws.Range("K31").Value = 1185.5454
$ws.Range("M31").Value = -890.5454
$ws.Range("H34").Value = 3469.4
$ws.Range("I34").Value = 1185.5454
$ws.Range("K34").Value = 1185.5454
$ws.Range("M34").Value = -983.5454
$ws.Range("H58").Value = 1739.1818
$ws.Range("I58").Value = 1558.2295
$ws.Range("K58").Value = 1558.2295
$ws.Range("M58").Value = -1355.2295
$ws.Range("H136").Value = 1739.1818
$ws.Range("I136").Value = 1558.2295
$ws.Range("K136").Value = 4674.6885
$ws.Range("M136").Value = -2124.6885
$ws.Range("H140").Value = 34776.668
$ws.Range("J140").Value = 34776.668
$ws.Range("L140").Value = 34776.668
$ws.Range("N140").Value = -45136.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 143
$ws.Range("I38").Value = 15
$ws.Range("J38").Value = 175
$ws.Range("K38").Value = 45
$ws.Range("L38").Value = 525
$ws.Range("M38").Value = 302
$ws.Range("N38").Value = -1219
$ws.Range("H113").Value = 669.72095
$ws.Range("I113").Value = 664.15625
$ws.Range("J113").Value = 685.9091
$ws.Range("K113").Value = 1992.46875
$ws.Range("L113").Value = 2057.7273
$ws.Range("M113").Value = 177.53125
$ws.Range("N113").Value = -6397.7273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6232.4287
$ws.Range("I70").Value = 5769.689
$ws.Range("J70").Value = 8125.4546
$ws.Range("K70").Value = 5769.689
$ws.Range("L70").Value = 8125.4546
$ws.Range("M70").Value = -5499.689
$ws.Range("N70").Value = -8665.454600000001
$ws.Range("H73").Value = 6232.4287
$ws.Range("I73").Value = 5769.689
$ws.Range("J73").Value = 8125.4546
$ws.Range("K73").Value = 5769.689
$ws.Range("L73").Value = 8125.4546
$ws.Range("M73").Value = -4833.689
$ws.Range("N73").Value = -9997.454600000001
$ws.Range("H97").Value = 785.6667
$ws.Range("I97").Value = 846.8570999999999
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 846.8570999999999
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = -350.8570999999999
$ws.Range("N97").Value = -1692
$ws.Range("H123").Value = 10493.125
$ws.Range("J123").Value = 10493.125
$ws.Range("L123").Value = 10493.125
$ws.Range("N123").Value = -15393.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 26900
$ws.Range("J64").Value = 26900
$ws.Range("L64").Value = 26900
$ws.Range("N64").Value = -27396
$ws.Range("H67").Value = 26900
$ws.Range("J67").Value = 26900
$ws.Range("L67").Value = 26900
$ws.Range("N67").Value = -28616
